$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new columns before column E (term) for female_learners / male_learners
$ws.Range("E1:F1").EntireColumn.Insert()

$ws.Range("E1").Value = "female_learners"
$ws.Range("F1").Value = "male_learners"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 4

# Column widths for the newly-shifted/added columns
# (target stored widths: D=10.83203125, E=15.6640625, F=17 -- set via the
# ColumnWidth property, which is the closest this host can reproduce them)
$ws.Range("D1").EntireColumn.ColumnWidth = 10
$ws.Range("E1").EntireColumn.ColumnWidth = 14.833333333333334
$ws.Range("F1").EntireColumn.ColumnWidth = 16.166666666666668

# Restore selection/view state
$ws.Range("G7").Select() | Out-Null
